$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Duplicate the formatting of the last existing data row (row 87) down into
# the two new rows, then overwrite with the new day's values.
$ws.Range("A87:F87").Copy()
$ws.Range("A88:F89").PasteSpecial(-4122)

$ws.Cells.Item(88, 1).Value = 45913
$ws.Cells.Item(88, 2).Value = "四方坪站"
$ws.Cells.Item(88, 3).Value = 11101.07
$ws.Cells.Item(88, 4).Value = 8986.36
$ws.Cells.Item(88, 5).Value = 3805.53
$ws.Cells.Item(88, 6).Value = 436

$ws.Cells.Item(89, 1).Value = 45913
$ws.Cells.Item(89, 2).Value = "高岭站"
$ws.Cells.Item(89, 3).Value = 7481.71
$ws.Cells.Item(89, 4).Value = 6233.88
$ws.Cells.Item(89, 5).Value = 1783.3
$ws.Cells.Item(89, 6).Value = 256

$ws.Range("H85").Select()
